# docs/ad_new_document_form.xlsx — "drew some guidline to create an order"
#
# Replaces the old scratch-work rows (12-27) that referenced internal
# pricing fields (product_cost, price_before_tax, tax, discount, ...) with
# a small "Necessary payload" guide table describing the objects needed to
# create an order (Sale / Document / Item / Payment, plus Nnumber / User /
# Customer / Cash register), and highlights the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- wipe the old (rows 12-27) scratch content ------------------------------
$ws.Range("C12:F27").ClearContents()

# --- new "necessary payload" guide ------------------------------------------
$ws.Range("C10").Value = "Necessary payload:"

$ws.Range("C12").Value = "Sale"
$ws.Range("D12").Value = "Document"
$ws.Range("E12").Value = "Item"
$ws.Range("F12").Value = "Payment"

$ws.Range("D13").Value = "Nnumber"
$ws.Range("D14").Value = "User"
$ws.Range("D15").Value = "Customer"
$ws.Range("D16").Value = "Cash register"

# highlight the header row (Sale / Document / Item / Payment) with the
# standard "Blue, Accent 1, Lighter 80%" theme fill
$ws.Range("C12:F12").Interior.Color = 15130576

# --- column C narrowed slightly ---------------------------------------------
$ws.Columns("C").ColumnWidth = 11.6

# --- selection cosmetically moved to just below the new table --------------
$ws.Range("D17").Select() | Out-Null
